$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (QSE) and rename it to CUSTUM
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "CUSTUM"

# Populate the new sheet with the custom list of stock symbols
$newSheet.Range("A1").Value = "SAZEW"
$newSheet.Range("A2").Value = "ICL"
$newSheet.Range("A3").Value = "IMAGE"
$newSheet.Range("A4").Value = "PNSC"

# Match the author's saved selection/active cell on the new sheet
$newSheet.Range("G16").Select()
